# "implemented HW1 for NHPA"
# Append four new log rows (17-20) to the time-tracking sheet:
#   - two sessions of "Numerical Algorithms" tutorium prep/STL work
#   - a new "STL" topic (two sessions)
#   - a new "Numerical HP Algorithms" (NHPA) homework session
# This introduces two brand-new shared strings ("STL" and
# "Numerical HP Algorithms") the same way Excel would: simply by typing
# them into cells that don't yet use those strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 is the last populated row and carries the correct per-column
# formatting (date format in A, time format in B:D). Clone that
# formatting down into the new rows before filling in values, mirroring
# how the original author kept extending the log by copying the row
# above. Row 20 never gets a date (A) value in the source data, so its
# formatting is copied separately, limited to B:E.
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E19").PasteSpecial(-4122)
$ws.Range("B16:E16").Copy()
$ws.Range("B20:E20").PasteSpecial(-4122)

# Row 17: 2018-09-27, 19:00-21:00, Numerical Algorithms, 2h
$ws.Range("A17").Value() = 43370
$ws.Range("B17").Value() = "Numerical Algorithms"
$ws.Range("C17").Value() = 0.79166666666666663
$ws.Range("D17").Value() = 0.875
$ws.Range("E17").Value() = 2

# Row 18: 2018-10-08, 09:00-11:00, STL, 2h
$ws.Range("A18").Value() = 43381
$ws.Range("B18").Value() = "STL"
$ws.Range("C18").Value() = 0.375
$ws.Range("D18").Value() = 0.45833333333333331
$ws.Range("E18").Value() = 2

# Row 19: 2018-10-15, 09:00-11:00, STL, 2h
$ws.Range("A19").Value() = 43388
$ws.Range("B19").Value() = "STL"
$ws.Range("C19").Value() = 0.375
$ws.Range("D19").Value() = 0.45833333333333331
$ws.Range("E19").Value() = 2

# Row 20: (no date) 16:30-17:30, Numerical HP Algorithms (HW1), 1h
$ws.Range("B20").Value() = "Numerical HP Algorithms"
$ws.Range("C20").Value() = 0.6875
$ws.Range("D20").Value() = 0.72916666666666663
$ws.Range("E20").Value() = 1

# Match the author's final cursor position/selection (cell below the
# newly entered data, ready for the next entry).
$ws.Range("A21").Select()
